$d = $word.ActiveDocument

$replacements = @(
    @{old="490÷5=98, 0"; new="204÷2=102, 0"},
    @{old="720÷8=90, 0"; new="249÷9=27, 6"},
    @{old="303÷6=50, 3"; new="341÷4=85, 1"},
    @{old="531÷7=75, 6"; new="515÷4=128, 3"},
    @{old="184÷8=23, 0"; new="104÷3=34, 2"},
    @{old="473÷7=67, 4"; new="678÷8=84, 6"},
    @{old="639÷5=127, 4"; new="565÷7=80, 5"},
    @{old="173÷3=57, 2"; new="308÷7=44, 0"},
    @{old="761÷6=126, 5"; new="759÷6=126, 3"},
    @{old="177÷7=25, 2"; new="923÷8=115, 3"},
    @{old="263÷2=131, 1"; new="259÷5=51, 4"},
    @{old="915÷5=183, 0"; new="916÷9=101, 7"},
    @{old="145÷6=24, 1"; new="338÷8=42, 2"},
    @{old="171÷2=85, 1"; new="248÷9=27, 5"},
    @{old="509÷3=169, 2"; new="897÷2=448, 1"},
    @{old="882÷9=98, 0"; new="627÷6=104, 3"},
    @{old="272÷7=38, 6"; new="521÷9=57, 8"},
    @{old="198÷5=39, 3"; new="827÷6=137, 5"},
    @{old="225÷6=37, 3"; new="274÷3=91, 1"},
    @{old="551÷6=91, 5"; new="574÷8=71, 6"},
    @{old="398÷6=66, 2"; new="891÷5=178, 1"},
    @{old="585÷6=97, 3"; new="750÷6=125, 0"},
    @{old="734÷9=81, 5"; new="929÷3=309, 2"},
    @{old="513÷9=57, 0"; new="234÷5=46, 4"},
    @{old="127÷7=18, 1"; new="260÷8=32, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
